# Applies the "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta" edit:
#  - Updates the total "VALOR MORA" figure (E11)
#  - Inserts one new period row (2508) at the end of the statement table
#  - Re-orders the period table from descending (2507 -> 1701) to ascending
#    (1701 -> 2508) chronological order, moving the "current" value (91917)
#    from the old first period (2507) to the new first period (1701) and
#    giving the new last period (2508) the old closing value (551500)
#  - Bumps the period counter (F13) to account for the new period

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update total overdue amount
$ws.Range("E11").Value = 56896417

# 2) Make room for the new period row: everything from row 119 down
#    (the blank gap + signature block) shifts down by one row, exactly
#    like Excel does when a new row is inserted above it.
$ws.Rows("119").Insert()

# 3) Row 119 should inherit the special "closing row" formatting that used
#    to live on row 118, and row 118 should become a normal data row like
#    the one above it (row 117).
$ws.Range("B118:J118").Copy($ws.Range("B119:J119"))
$ws.Range("B117:J117").Copy($ws.Range("B118:J118"))

# 4) Build the chronological list of periods: 1701 .. 2506, then 2507, 2508
$periods = @()
for ($yy = 17; $yy -le 24; $yy++) {
    for ($mm = 1; $mm -le 12; $mm++) {
        $periods += ("{0:D2}{1:D2}" -f $yy, $mm)
    }
}
for ($mm = 1; $mm -le 6; $mm++) {
    $periods += ("25{0:D2}" -f $mm)
}
$periods += "2507"
$periods += "2508"

# 5) Write the period label into each row (B16:J119), and the "Valor Mora"
#    amount: the first (oldest) period carries the 91917 figure, every
#    other period keeps the standard 551500 figure.
$firstRow = 16
$lastRow = 119
for ($i = 0; $i -lt $periods.Count; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    if ($i -eq 0) {
        $ws.Cells.Item($row, 6).Value = 91917
    } else {
        $ws.Cells.Item($row, 6).Value = 551500
    }
}

# 6) Bump the "Cant. Periodos" counter to reflect the newly added period
$ws.Range("F13").Value = 104
